$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings that often look numeric (e.g. "1.006", "0.00001067").
# Force text formatting first so Excel stores them verbatim instead of
# reinterpreting them as numbers/dates.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '27.594.39'
$ws.Range("E2").Value = '  -2.58%  '
$ws.Range("D3").Value = '1.750.10'
$ws.Range("E3").Value = '  -3.74%  '
$ws.Range("D4").Value = '1.006'
$ws.Range("E4").Value = '  +0.53%  '
$ws.Range("D5").Value = '321.95'
$ws.Range("E5").Value = '  -2.30%  '
$ws.Range("D6").Value = '1.004'
$ws.Range("E6").Value = '  +0.46%  '
$ws.Range("D7").Value = '0.4235'
$ws.Range("E7").Value = '  -4.28%  '
$ws.Range("D8").Value = '0.3616'
$ws.Range("E8").Value = '  -2.45%  '
$ws.Range("B9").Value = 'OKB'
$ws.Range("C9").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D9").Value = '42.50'
$ws.Range("E9").Value = '  -4.72%  '
$ws.Range("B10").Value = 'Dogecoin'
$ws.Range("C10").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D10").Value = '0.07491'
$ws.Range("E10").Value = '  -2.90%  '
$ws.Range("D11").Value = '1.090'
$ws.Range("E11").Value = '  -3.46%  '
$ws.Range("D12").Value = '1.005'
$ws.Range("E12").Value = '  +0.63%  '
$ws.Range("D13").Value = '20.50'
$ws.Range("E13").Value = '  -7.30%  '
$ws.Range("D14").Value = '6.004'
$ws.Range("E14").Value = '  -4.20%  '
$ws.Range("D15").Value = '7.235'
$ws.Range("E15").Value = '  -4.18%  '
$ws.Range("D16").Value = '1.773.58'
$ws.Range("E16").Value = '  -2.41%  '
$ws.Range("D17").Value = '90.70'
$ws.Range("E17").Value = '  -2.39%  '
$ws.Range("D18").Value = '0.00001067'
$ws.Range("E18").Value = '  -1.60%  '
$ws.Range("D19").Value = '0.06355'
$ws.Range("E19").Value = '  -4.55%  '
$ws.Range("D20").Value = '1.004'
$ws.Range("E20").Value = '  +0.45%  '
$ws.Range("D21").Value = '16.95'
$ws.Range("E21").Value = '  -3.64%  '
$ws.Range("D22").Value = '5.866'
$ws.Range("E22").Value = '  -5.79%  '
$ws.Range("D23").Value = '27.668.87'
$ws.Range("E23").Value = '  -2.46%  '
$ws.Range("D24").Value = '11.12'
$ws.Range("E24").Value = '  -5.01%  '
$ws.Range("D25").Value = '2.091'
$ws.Range("E25").Value = '  +5.30%  '
$ws.Range("D26").Value = '160.61'
$ws.Range("E26").Value = '  +2.84%  '
$ws.Range("D27").Value = '20.15'
$ws.Range("E27").Value = '  -2.97%  '
$ws.Range("D28").Value = '1.979.77'
$ws.Range("E28").Value = '  -2.15%  '
$ws.Range("D29").Value = '2.120'
$ws.Range("E29").Value = '  -8.67%  '
$ws.Range("D30").Value = '124.23'
$ws.Range("E30").Value = '  -3.20%  '
$ws.Range("D31").Value = '1.094'
$ws.Range("E31").Value = '  -9.35%  '
$ws.Range("D32").Value = '3.656'
$ws.Range("E32").Value = '  -0.23%  '
$ws.Range("D33").Value = '5.529'
$ws.Range("E33").Value = '  -5.91%  '
$ws.Range("D34").Value = '0.08863'
$ws.Range("E34").Value = '  -3.91%  '
$ws.Range("D35").Value = '12.15'
$ws.Range("E35").Value = '  -7.12%  '
$ws.Range("D36").Value = '0.02280'
$ws.Range("E36").Value = '  -3.52%  '
$ws.Range("D37").Value = '0.2094'
$ws.Range("E37").Value = '  -3.75%  '
$ws.Range("D38").Value = '0.06009'
$ws.Range("E38").Value = '  -3.49%  '
$ws.Range("D39").Value = '0.6310'
$ws.Range("E39").Value = '  -4.13%  '
$ws.Range("D40").Value = '4.930'
$ws.Range("E40").Value = '  -4.64%  '
$ws.Range("D41").Value = '1.182'
$ws.Range("E41").Value = '  -1.49%  '
$ws.Range("D42").Value = '1.004'
$ws.Range("E42").Value = '  +0.48%  '
$ws.Range("D43").Value = '7.848'
$ws.Range("E43").Value = '  -3.74%  '
$ws.Range("D44").Value = '1.395'
$ws.Range("E44").Value = '  +0.58%  '
$ws.Range("D45").Value = '13.23'
$ws.Range("E45").Value = '  -4.82%  '
$ws.Range("D46").Value = '0.5852'
$ws.Range("E46").Value = '  -3.97%  '
$ws.Range("D47").Value = '3.689'
$ws.Range("E47").Value = '  -2.00%  '
$ws.Range("B48").Value = 'Quant'
$ws.Range("C48").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D48").Value = '123.15'
$ws.Range("E48").Value = '  -2.83%  '
$ws.Range("B49").Value = 'NEARProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D49").Value = '1.977'
$ws.Range("E49").Value = '  -3.08%  '
$ws.Range("D50").Value = '1.165'
$ws.Range("E50").Value = '  +0.87%  '
$ws.Range("D51").Value = '0.06811'
$ws.Range("E51").Value = '  -2.46%  '
Write-Host "Updated cryptos list"
